# Update scripts with new TPM values: refresh the ligand/receptor expression
# and edge-weight statistics (columns G,H,I,J,M,N,O,P,Q,R,S,T) for rows 2-10
# to reflect the recomputed TPM-based figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.6244246666666666
$ws.Cells.Item(2, 8).Value = 1.873274
$ws.Cells.Item(2, 9).Value = 0.1217547960993207
$ws.Cells.Item(2, 10).Value = 0.1217547960993207
$ws.Cells.Item(2, 13).Value = 4.685980666666667
$ws.Cells.Item(2, 14).Value = 14.057942
$ws.Cells.Item(2, 15).Value = 0.2075482316885617
$ws.Cells.Item(2, 16).Value = 0.2075482316885617
$ws.Cells.Item(2, 17).Value = 2.926041915789777
$ws.Cells.Item(2, 18).Value = 26.334377242108
$ws.Cells.Item(2, 19).Value = 0.02526999263001541
$ws.Cells.Item(2, 20).Value = 0.02526999263001541

$ws.Cells.Item(3, 7).Value = 0.6244246666666666
$ws.Cells.Item(3, 8).Value = 1.873274
$ws.Cells.Item(3, 9).Value = 0.1217547960993207
$ws.Cells.Item(3, 10).Value = 0.1217547960993207
$ws.Cells.Item(3, 15).Value = 0.5625283482851452
$ws.Cells.Item(3, 16).Value = 0.5625283482851452
$ws.Cells.Item(3, 17).Value = 7.930597685709111
$ws.Cells.Item(3, 18).Value = 71.375379171382
$ws.Cells.Item(3, 19).Value = 0.06849052434554553
$ws.Cells.Item(3, 20).Value = 0.06849052434554553

$ws.Cells.Item(4, 7).Value = 0.6244246666666666
$ws.Cells.Item(4, 8).Value = 1.873274
$ws.Cells.Item(4, 9).Value = 0.1217547960993207
$ws.Cells.Item(4, 10).Value = 0.1217547960993207
$ws.Cells.Item(4, 13).Value = 5.191163
$ws.Cells.Item(4, 14).Value = 15.573489
$ws.Cells.Item(4, 15).Value = 0.2299234200262931
$ws.Cells.Item(4, 16).Value = 0.2299234200262931
$ws.Cells.Item(4, 17).Value = 3.241490225887333
$ws.Cells.Item(4, 18).Value = 29.173412032986
$ws.Cells.Item(4, 19).Value = 0.02799427912375979
$ws.Cells.Item(4, 20).Value = 0.02799427912375979

$ws.Cells.Item(5, 7).Value = 3.569858333333333
$ws.Cells.Item(5, 9).Value = 0.6960765592408708
$ws.Cells.Item(5, 10).Value = 0.6960765592408706
$ws.Cells.Item(5, 13).Value = 4.685980666666667
$ws.Cells.Item(5, 14).Value = 14.057942
$ws.Cells.Item(5, 15).Value = 0.2075482316885617
$ws.Cells.Item(5, 16).Value = 0.2075482316885617
$ws.Cells.Item(5, 17).Value = 16.72828713273889
$ws.Cells.Item(5, 18).Value = 150.55458419465
$ws.Cells.Item(5, 19).Value = 0.1444694589903011
$ws.Cells.Item(5, 20).Value = 0.1444694589903011

$ws.Cells.Item(6, 7).Value = 3.569858333333333
$ws.Cells.Item(6, 9).Value = 0.6960765592408708
$ws.Cells.Item(6, 10).Value = 0.6960765592408706
$ws.Cells.Item(6, 15).Value = 0.5625283482851452
$ws.Cells.Item(6, 16).Value = 0.5625283482851452
$ws.Cells.Item(6, 17).Value = 45.33951291158055
$ws.Cells.Item(6, 19).Value = 0.391562797149774
$ws.Cells.Item(6, 20).Value = 0.391562797149774

$ws.Cells.Item(7, 7).Value = 3.569858333333333
$ws.Cells.Item(7, 9).Value = 0.6960765592408708
$ws.Cells.Item(7, 10).Value = 0.6960765592408706
$ws.Cells.Item(7, 13).Value = 5.191163
$ws.Cells.Item(7, 14).Value = 15.573489
$ws.Cells.Item(7, 15).Value = 0.2299234200262931
$ws.Cells.Item(7, 16).Value = 0.2299234200262931
$ws.Cells.Item(7, 17).Value = 18.53171649524166
$ws.Cells.Item(7, 18).Value = 166.785448457175
$ws.Cells.Item(7, 19).Value = 0.1600443031007956
$ws.Cells.Item(7, 20).Value = 0.1600443031007956

$ws.Cells.Item(8, 7).Value = 0.9342596666666667
$ws.Cells.Item(8, 8).Value = 2.802779
$ws.Cells.Item(8, 9).Value = 0.1821686446598085
$ws.Cells.Item(8, 10).Value = 0.1821686446598085
$ws.Cells.Item(8, 13).Value = 4.685980666666667
$ws.Cells.Item(8, 14).Value = 14.057942
$ws.Cells.Item(8, 15).Value = 0.2075482316885617
$ws.Cells.Item(8, 16).Value = 0.2075482316885617
$ws.Cells.Item(8, 17).Value = 4.377922735646445
$ws.Cells.Item(8, 18).Value = 39.401304620818
$ws.Cells.Item(8, 19).Value = 0.0378087800682452
$ws.Cells.Item(8, 20).Value = 0.03780878006824521

$ws.Cells.Item(9, 7).Value = 0.9342596666666667
$ws.Cells.Item(9, 8).Value = 2.802779
$ws.Cells.Item(9, 9).Value = 0.1821686446598085
$ws.Cells.Item(9, 10).Value = 0.1821686446598085
$ws.Cells.Item(9, 15).Value = 0.5625283482851452
$ws.Cells.Item(9, 16).Value = 0.5625283482851452
$ws.Cells.Item(9, 17).Value = 11.86570285551078
$ws.Cells.Item(9, 18).Value = 106.791325699597
$ws.Cells.Item(9, 19).Value = 0.1024750267898256
$ws.Cells.Item(9, 20).Value = 0.1024750267898256

$ws.Cells.Item(10, 7).Value = 0.9342596666666667
$ws.Cells.Item(10, 8).Value = 2.802779
$ws.Cells.Item(10, 9).Value = 0.1821686446598085
$ws.Cells.Item(10, 10).Value = 0.1821686446598085
$ws.Cells.Item(10, 13).Value = 5.191163
$ws.Cells.Item(10, 14).Value = 15.573489
$ws.Cells.Item(10, 15).Value = 0.2299234200262931
$ws.Cells.Item(10, 16).Value = 0.2299234200262931
$ws.Cells.Item(10, 17).Value = 4.849894213992333
$ws.Cells.Item(10, 18).Value = 43.649047925931
$ws.Cells.Item(10, 19).Value = 0.04188483780173768
$ws.Cells.Item(10, 20).Value = 0.04188483780173768
